# Penalty Reward System (unfinished) - shift all Week_Start_Date values
# forward by one week (7 days) across the "Forecast Comparison" sheet and
# update the related summary figures on the "Summary" sheet.
#
# Dates are stored as plain text (yyyy-mm-dd), not as real Excel date
# serials, so each cell is forced to Text format ("@") before writing the
# new literal (otherwise Excel auto-parses a string like "2025-01-12"
# into a date serial number). The format is reset back to Normal
# afterwards so no stray formatting is left behind on the cell.
# (Note: COM range objects must be accessed directly / via script-scope
# variables here rather than passed through function parameters, which
# this runtime does not marshal correctly.)

$wb = $excel.ActiveWorkbook

# ----- Sheet 1: Forecast Comparison -----
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")

$newDates = @(
    "2025-01-12",
    "2025-01-19",
    "2025-01-26",
    "2025-02-02",
    "2025-02-09",
    "2025-02-16",
    "2025-02-23",
    "2025-03-02",
    "2025-03-09",
    "2025-03-16",
    "2025-03-23",
    "2025-03-30",
    "2025-04-06",
    "2025-04-13",
    "2025-04-20",
    "2025-04-27"
)

for ($i = 0; $i -lt $newDates.Length; $i++) {
    $row = $i + 2
    $wsForecast.Cells.Item($row, 2).NumberFormat = "@"
    $wsForecast.Cells.Item($row, 2).Value2 = $newDates[$i]
    $wsForecast.Cells.Item($row, 2).Style = "Normal"
}

# ----- Sheet 2: Summary -----
$wsSummary = $wb.Worksheets.Item("Summary")

$wsSummary.Range("B2").NumberFormat = "@"
$wsSummary.Range("B2").Value2 = "2023-01-01 to 2025-01-05"
$wsSummary.Range("B2").Style = "Normal"

$wsSummary.Range("B13").NumberFormat = "@"
$wsSummary.Range("B13").Value2 = "2025-01-12"
$wsSummary.Range("B13").Style = "Normal"

$wsSummary.Range("B15").NumberFormat = "@"
$wsSummary.Range("B15").Value2 = "2025-01-12"
$wsSummary.Range("B15").Style = "Normal"
